# Supported.xlsx update: mark several additional BASIC functions as
# supported (columns B, C and D get an "X") and move the active
# selection down to where the most recent edits were made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose B/C/D "supported" columns need to be filled in with "X"
$rows = @(13, 23, 24, 26, 27, 28, 51, 96, 97, 98)

foreach ($r in $rows) {
    $ws.Range("B$r").Value2 = "X"
    $ws.Range("C$r").Value2 = "X"
    $ws.Range("D$r").Value2 = "X"
}

# Scroll the view down near the bottom of the edits and leave the
# selection on D51, matching where the author ended up working.
try {
    $excel.ActiveWindow.ScrollRow = 45
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Older/limited hosts may not support programmatic scrolling;
    # the selection change below is the important part.
}
$ws.Range("D51").Select()
